$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '69.391.70'
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  +0.07%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '3.672.81'
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('E4').Value = '  +0.00%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '645.16'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -5.23%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '159.72'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('E9').Value = '  -0.43%  '
$ws.Range('E10').Value = '  -0.52%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.449'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  +2.06%  '
$ws.Range('E12').Value = '  +0.58%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '4.297.28'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  -0.39%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '32.77'
$cell.Style = 'Normal'
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '3.673.19'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -0.33%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '69.416.37'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('E17').Value = '  +0.21%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '16.01'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  +0.04%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '6.50'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  +0.09%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '466.00'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  -0.58%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '9.91'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  -0.04%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '0.646'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  -1.25%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '79.46'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  -0.66%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '3.820.60'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  -0.44%  '
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('E26').Value = '  +2.79%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '10.89'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -0.27%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '9.05'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  -0.93%  '
$ws.Range('E29').Value = '  -2.53%  '
$ws.Range('E30').Value = '  -1.05%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '2.01'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  +0.49%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  -0.22%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '26.86'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -0.34%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '0.165'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  +4.62%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '6.46'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -1.86%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '3.666.56'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  -0.33%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '8.42'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  +1.22%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '5.89'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  -5.56%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '178.39'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  +4.14%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  +0.03%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '0.0899'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  -0.50%  '
$ws.Range('E43').Value = '  -1.71%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '0.925'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  -1.88%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '46.63'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  -1.89%  '
$ws.Range('E46').Value = '  +1.58%  '
$ws.Range('E47').Value = '  -1.81%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '27.20'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  -3.71%  '
$ws.Range('E49').Value = '  -4.78%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '1.25'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -4.02%  '
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '7.84'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  +0.56%  '
